$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-13 Thursday", "2025-03-14 Friday"),
    @("37÷7=5, 2", "73÷6=12, 1"),
    @("15÷8=1, 7", "94÷7=13, 3"),
    @("27÷7=3, 6", "68÷6=11, 2"),
    @("81÷2=40, 1", "45÷8=5, 5"),
    @("82÷3=27, 1", "33÷6=5, 3"),
    @("98÷5=19, 3", "87÷6=14, 3"),
    @("83÷7=11, 6", "91÷3=30, 1"),
    @("45÷7=6, 3", "81÷4=20, 1"),
    @("28÷7=4, 0", "88÷8=11, 0"),
    @("18÷8=2, 2", "96÷7=13, 5"),
    @("22÷7=3, 1", "59÷9=6, 5"),
    @("74÷8=9, 2", "88÷2=44, 0"),
    @("24÷4=6, 0", "33÷9=3, 6"),
    @("23÷6=3, 5", "54÷2=27, 0"),
    @("23÷7=3, 2", "34÷4=8, 2"),
    @("59÷3=19, 2", "35÷6=5, 5"),
    @("46÷9=5, 1", "63÷3=21, 0"),
    @("74÷4=18, 2", "96÷6=16, 0"),
    @("24÷3=8, 0", "35÷9=3, 8"),
    @("70÷8=8, 6", "84÷3=28, 0"),
    @("84÷4=21, 0", "54÷6=9, 0"),
    @("34÷2=17, 0", "61÷8=7, 5"),
    @("42÷4=10, 2", "89÷4=22, 1"),
    @("67÷7=9, 4", "69÷6=11, 3"),
    @("43÷2=21, 1", "93÷3=31, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
